$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 5267514
$ws.Range("I88").Value = 7844.3335
$ws.Range("J88").Value = 10001217
$ws.Range("K88").Value = 7844.3335
$ws.Range("L88").Value = 10001217
$ws.Range("M88").Value = -7438.3335
$ws.Range("N88").Value = -10002029

$ws.Range("H91").Value = 5267514
$ws.Range("I91").Value = 7844.3335
$ws.Range("J91").Value = 10001217
$ws.Range("K91").Value = 7844.3335
$ws.Range("L91").Value = 10001217
$ws.Range("M91").Value = -6440.3335
$ws.Range("N91").Value = -10004025

$ws.Range("H111").Value = 2102.1052
$ws.Range("J111").Value = 1384.7273
$ws.Range("L111").Value = 4154.1819
$ws.Range("N111").Value = -10288.1819

$ws.Range("H113").Value = 160856.78
$ws.Range("I113").Value = 287526.2
$ws.Range("J113").Value = 2520
$ws.Range("K113").Value = 287526.2
$ws.Range("L113").Value = 2520
$ws.Range("M113").Value = -284272.2
$ws.Range("N113").Value = -9028

$ws.Range("H132").Value = 5688003.5
$ws.Range("I132").Value = 5471.1
$ws.Range("J132").Value = 17864858
$ws.Range("K132").Value = 16413.3
$ws.Range("L132").Value = 53594574
$ws.Range("M132").Value = -13883.3
$ws.Range("N132").Value = -53599634

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1970.3334
$ws.Range("I2").Value = 1755.5
$ws.Range("J2").Value = 2400
$ws.Range("K2").Value = 1755.5
$ws.Range("L2").Value = 2400
$ws.Range("M2").Value = -1642.5
$ws.Range("N2").Value = -2626

$ws.Range("H32").Value = 6681.7744
$ws.Range("I32").Value = 7023.9805
$ws.Range("K32").Value = 7023.9805
$ws.Range("M32").Value = -6736.9805

$ws.Range("H45").Value = 5971.4287
$ws.Range("I45").Value = 6350
$ws.Range("J45").Value = 5466.6665
$ws.Range("K45").Value = 6350
$ws.Range("L45").Value = 5466.6665
$ws.Range("M45").Value = -5973
$ws.Range("N45").Value = -6220.6665

$ws.Range("H61").Value = 12502413
$ws.Range("I61").Value = 13891156
$ws.Range("J61").Value = 3724.5
$ws.Range("K61").Value = 13891156
$ws.Range("L61").Value = 3724.5
$ws.Range("M61").Value = -13890944
$ws.Range("N61").Value = -4148.5

$ws.Range("H116").Value = 1970.3334
$ws.Range("I116").Value = 1755.5
$ws.Range("J116").Value = 2400
$ws.Range("K116").Value = 1755.5
$ws.Range("L116").Value = 2400
$ws.Range("M116").Value = 538.5
$ws.Range("N116").Value = -6988

$ws.Range("H132").Value = 6946616
$ws.Range("I132").Value = 11906606
$ws.Range("J132").Value = 2629
$ws.Range("K132").Value = 35719818
$ws.Range("L132").Value = 7887
$ws.Range("M132").Value = -35717288
$ws.Range("N132").Value = -12947

$ws.Range("H136").Value = 12502413
$ws.Range("I136").Value = 13891156
$ws.Range("J136").Value = 3724.5
$ws.Range("K136").Value = 41673468
$ws.Range("L136").Value = 11173.5
$ws.Range("M136").Value = -41670918
$ws.Range("N136").Value = -16273.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1970.3334
$ws.Range("I3").Value = 1755.5
$ws.Range("J3").Value = 2400
$ws.Range("K3").Value = 1755.5
$ws.Range("L3").Value = 2400
$ws.Range("M3").Value = -1641.5
$ws.Range("N3").Value = -2628

$ws.Range("H107").Value = 2052.6924
$ws.Range("I107").Value = 2228.45
$ws.Range("J107").Value = 1466.8334
$ws.Range("K107").Value = 2228.45
$ws.Range("L107").Value = 1466.8334
$ws.Range("M107").Value = -308.4499999999998
$ws.Range("N107").Value = -5306.8334

$ws.Range("H134").Value = 4009.75
$ws.Range("I134").Value = 2764.6897
$ws.Range("K134").Value = 8294.069100000001
$ws.Range("M134").Value = -5759.069100000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1187.375
$ws.Range("I16").Value = 999.8
$ws.Range("J16").Value = 1500
$ws.Range("K16").Value = 999.8
$ws.Range("L16").Value = 1500
$ws.Range("M16").Value = -712.8
$ws.Range("N16").Value = -2074

$ws.Range("H58").Value = 2493.8928
$ws.Range("I58").Value = 1826.4615
$ws.Range("K58").Value = 1826.4615
$ws.Range("M58").Value = -1623.4615

$ws.Range("H113").Value = 1187.375
$ws.Range("I113").Value = 999.8
$ws.Range("J113").Value = 1500
$ws.Range("K113").Value = 999.8
$ws.Range("L113").Value = 1500
$ws.Range("M113").Value = 1170.2
$ws.Range("N113").Value = -5840

$ws.Range("H132").Value = 4954.1816
$ws.Range("I132").Value = 5271.2
$ws.Range("K132").Value = 15813.6
$ws.Range("M132").Value = -13283.6

$ws.Range("H134").Value = 955329.2
$ws.Range("I134").Value = 3582
$ws.Range("K134").Value = 10746
$ws.Range("M134").Value = -8211

$ws.Range("H136").Value = 2493.8928
$ws.Range("I136").Value = 1826.4615
$ws.Range("K136").Value = 5479.3845
$ws.Range("M136").Value = -2929.3845

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 5585.35
$ws.Range("I3").Value = 3943.182
$ws.Range("J3").Value = 7592.4443
$ws.Range("K3").Value = 11829.546
$ws.Range("L3").Value = 22777.3329
$ws.Range("M3").Value = -11717.546
$ws.Range("N3").Value = -23001.3329

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 5448.3105
$ws.Range("I132").Value = 5151.467
$ws.Range("J132").Value = 5766.357
$ws.Range("K132").Value = 15454.401
$ws.Range("L132").Value = 17299.071
$ws.Range("M132").Value = -12924.401
$ws.Range("N132").Value = -22359.071

$ws.Range("H134").Value = 26751.273
$ws.Range("J134").Value = 26751.273
$ws.Range("L134").Value = 80253.819
$ws.Range("N134").Value = -85323.819

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 3074.4167
$ws.Range("I100").Value = 4082.1667
$ws.Range("K100").Value = 8164.3334
$ws.Range("M100").Value = -7623.3334

$ws.Range("H107").Value = 1903.4
$ws.Range("I107").Value = 2329.1428
$ws.Range("J107").Value = 910
$ws.Range("K107").Value = 6987.428400000001
$ws.Range("L107").Value = 2730
$ws.Range("M107").Value = -5067.428400000001
$ws.Range("N107").Value = -6570

$ws.Range("H137").Value = 67567.39999999999
$ws.Range("J137").Value = 67567.39999999999
$ws.Range("L137").Value = 67567.39999999999
$ws.Range("N137").Value = -77767.39999999999
